$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Numeric value updates ---
$ws.Range("H3").Value = 0.0460028648376464
$ws.Range("I3").Value = 0.0560040473937988

$ws.Range("I7").Value = 0.00500106811523437

# H15 / I15 values updated (new H15 value, old H15 value moved into I15)
$ws.Range("H15").Value = 0.0160009860992431
$ws.Range("I15").Value = 0.0410029888153076

$ws.Range("H19").Value = 0.00399994850158691
$ws.Range("I19").Value = 0.00500082969665527

$ws.Range("H25").Value = 0.00600099563598632
$ws.Range("I25").Value = 0.00499987602233886

$ws.Range("H26").Value = 0.00500106811523437
$ws.Range("I26").Value = 0.00500011444091796

$ws.Range("I27").Value = 0.00500106811523437

$ws.Range("H28").Value = 0.00600099563598632
$ws.Range("I28").Value = 0.00600004196166992

# --- Shared string (text) value updates ---
$ws.Range("H9").Value = "~9.2x"
$ws.Range("I9").Value = "~11.2x"

$ws.Range("H21").Value = "~4x"
$ws.Range("I21").Value = "8.2x"
